# 8.5.2.xlsx — add a new "2022" column (S) mirroring the existing "2021"
# column (R): same per-row formatting, new data values, plus a fresh
# bold+italic style used on the two blank section-separator rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value for column S (mirrors column R one column over).
# Rows 8 and 36 are blank separator rows and are handled separately below.
$values = @(
    @(4,  2022),
    @(5,  4.9000000000000004),
    @(6,  6.1),
    @(7,  4),
    @(9,  6.1),
    @(10, 12.4),
    @(11, 3.2),
    @(12, 10.8),
    @(13, 14.6),
    @(14, 8.5),
    @(15, 5.5),
    @(16, 7.1),
    @(17, 4.4000000000000004),
    @(18, 5.8),
    @(19, 11.6),
    @(20, 3.1),
    @(21, 1.5),
    @(22, 2.2999999999999998),
    @(23, 1),
    @(24, 2.2999999999999998),
    @(25, 3.3),
    @(26, 1.6),
    @(27, 4.5999999999999996),
    @(28, 4.4000000000000004),
    @(29, 4.7),
    @(30, 4),
    @(31, 3.2),
    @(32, 4.7),
    @(33, 2.6),
    @(34, 3.3),
    @(35, 2.2000000000000002),
    @(37, 13.2),
    @(38, 7.5),
    @(39, 4.0999999999999996),
    @(40, 4.3),
    @(41, 2.6),
    @(42, 1)
)

foreach ($pair in $values) {
    $row = $pair[0]
    $val = $pair[1]

    # Copy column R's formatting for this row onto column S (reuses the
    # existing cellXf/font instead of minting new ones), then write the value.
    $src = $ws.Range("R" + $row)
    $src.Copy()
    $dst = $ws.Range("S" + $row)
    $dst.PasteSpecial(-4122)
    $dst.Value = $val
}
$excel.CutCopyMode = $false

# Rows 8 and 36 stay blank in column S, but pick up a brand new bold+italic
# 9pt Times New Roman style (copy the plain style from R, then embolden it).
foreach ($row in @(8, 36)) {
    $src = $ws.Range("R" + $row)
    $src.Copy()
    $dst = $ws.Range("S" + $row)
    $dst.PasteSpecial(-4122)
    $dst.Font.Bold = $true
    $dst.Font.Italic = $true
}
$excel.CutCopyMode = $false

# Row 43's "…" placeholder also carries over, reusing the shared string.
$src = $ws.Range("R43")
$src.Copy()
$dst = $ws.Range("S43")
$dst.PasteSpecial(-4122)
$dst.Value = "…"
$excel.CutCopyMode = $false

# Match the author's final selection/cursor position.
$null = $ws.Range("T12").Select()
